# Adding 6 new test cases to Search Module (TestCase_B35 .. TestCase_B40)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 36 : TestCase_B35 -------------------------------------------------
$ws.Cells.Item(36, 1).Value = "TestCase_B35"
$ws.Cells.Item(36, 2).Value = "TBD-10"
$ws.Cells.Item(36, 3).Value = "Verify that no filtering options are present in ALL content type"
$ws.Cells.Item(36, 4).Value = "Y"
$ws.Cells.Item(36, 5).Value = "PASS"

# --- Row 37 : TestCase_B36 (wrap-text description) -------------------------
$ws.Cells.Item(37, 1).Value = "TestCase_B36"
$ws.Cells.Item(37, 2).Value = "TBD-11"
$ws.Cells.Item(37, 3).Value = "Verify that the following fields get displayed in the SORT BY drop down when ARTICLES is selected as content type in the left navigation pane:`na)Relevance`nb)Times cited`nc)Publication Date(Newest)`nd)Publication Date(Oldest)`n"
$ws.Cells.Item(37, 4).Value = "Y"
$ws.Cells.Item(37, 5).Value = "PASS"

# --- Row 38 : TestCase_B37 --------------------------------------------------
$ws.Cells.Item(38, 1).Value = "TestCase_B37"
$ws.Cells.Item(38, 2).Value = "TBD-12"
$ws.Cells.Item(38, 3).Value = "Verify that user is able to sort the articles by TIMES CITED field in ARTICLES content type"
$ws.Cells.Item(38, 4).Value = "Y"
$ws.Cells.Item(38, 5).Value = "PASS"

# --- Row 39 : TestCase_B38 --------------------------------------------------
$ws.Cells.Item(39, 1).Value = "TestCase_B38"
$ws.Cells.Item(39, 2).Value = "TBD-13"
$ws.Cells.Item(39, 3).Value = "Verify that only articles get displayed when user chooses ARTICLES as content type"
$ws.Cells.Item(39, 4).Value = "Y"
$ws.Cells.Item(39, 5).Value = "PASS"

# --- Row 40 : TestCase_B39 --------------------------------------------------
$ws.Cells.Item(40, 1).Value = "TestCase_B39"
$ws.Cells.Item(40, 2).Value = "TBD-14"
$ws.Cells.Item(40, 3).Value = "Verify that all articles are sorted by RELEVANCE by default in ARTICLES content type"
$ws.Cells.Item(40, 4).Value = "Y"
$ws.Cells.Item(40, 5).Value = "PASS"

# --- Row 41 : TestCase_B40 (wrap-text description) --------------------------
$ws.Cells.Item(41, 1).Value = "TestCase_B40"
$ws.Cells.Item(41, 2).Value = "TBD-15"
$ws.Cells.Item(41, 3).Value = "Verify that following filters are present for ARTICLES content type:`na)Document Type`nb)Authors`nc)Categories`nd)Institutions`n"
$ws.Cells.Item(41, 4).Value = "Y"
$ws.Cells.Item(41, 5).Value = "PASS"

# --- Formatting: copy the existing bordered styles onto the new rows -------
# Plain bordered style (style index 3), the same one already used on A2, is
# applied to the whole new block first ...
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A36:E41").PasteSpecial(-4122) | Out-Null

# ... then the bordered + wrap-text style (style index 4) used by the long,
# multi-line descriptions is applied on top of C37 and C41 (same style as
# the existing C27 cell already uses).
$ws.Range("C27").Copy() | Out-Null
$ws.Range("C37").PasteSpecial(-4122) | Out-Null
$ws.Range("C41").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- View state: scroll so row 25 is at the top and select the new last cell
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A41").Select() | Out-Null
